$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 72, shifting existing rows 72-78 down to 73-79
$ws.Rows.Item(72).Insert()

# Fill in the new row 72 data
$ws.Cells.Item(72, 1).Value = 7
$ws.Cells.Item(72, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(72, 3).Value = "Ñuble"
$ws.Cells.Item(72, 4).Value = 45132
$ws.Cells.Item(72, 5).Value = 16
$ws.Cells.Item(72, 6).Value = 100112026
$ws.Cells.Item(72, 7).Value = "Haba"
$ws.Cells.Item(72, 8).Value = "Sin especificar"
$ws.Cells.Item(72, 9).Value = "Primera"
$ws.Cells.Item(72, 10).Value = 25
$ws.Cells.Item(72, 11).Value = 16000
$ws.Cells.Item(72, 12).Value = 16000
$ws.Cells.Item(72, 13).Value = 16000
$ws.Cells.Item(72, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(72, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(72, 16).Value = 640
$ws.Cells.Item(72, 17).Value = 25
$ws.Cells.Item(72, 18).Value = "Hortaliza"
